# Romania Liga I - atualizacao de bases (29-03-2024 17:05)
#
# The source feed re-fetched odds for three already-settled fixtures
# (rows 235-237, id 233-235) which came back in a different order, and
# the upcoming-fixtures tail (rows 251-256, id 249-254) shifted up by
# one row because one of the placeholder fixtures (id 249 / match
# 7951750) dropped out of the feed entirely.
#
# Net effect on the worksheet:
#   - row 235 (id 233) <- old row 237 data
#   - row 236 (id 234) <- old row 235 data
#   - row 237 (id 235) <- old row 236 data
#   - row 251 (id 249) <- old row 252 data
#   - row 252 (id 250) <- old row 253 data
#   - row 253 (id 251) <- old row 254 data
#   - row 254 (id 252) <- old row 255 data
#   - row 255 (id 253) <- old row 256 data
#   - row 256 removed (no longer present in the feed)
#
# Column A (the running id) is positional and must stay put; every
# other column (B:AC) moves with its row's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function RowRange($row) {
    return $ws.Range("B" + $row + ":AC" + $row)
}

# --- capture the "before" state of every row involved, before writing anything ---
$row235 = (RowRange 235).Value2
$row236 = (RowRange 236).Value2
$row237 = (RowRange 237).Value2

$row252 = (RowRange 252).Value2
$row253 = (RowRange 253).Value2
$row254 = (RowRange 254).Value2
$row255 = (RowRange 255).Value2
$row256 = (RowRange 256).Value2

# --- rotate rows 235-237 ---
(RowRange 235).Value = $row237
(RowRange 236).Value = $row235
(RowRange 237).Value = $row236

# --- shift rows 251-255 up from 252-256 ---
(RowRange 251).Value = $row252
(RowRange 252).Value = $row253
(RowRange 253).Value = $row254
(RowRange 254).Value = $row255
(RowRange 255).Value = $row256

# --- the old last row (256) no longer exists in the feed ---
$ws.Rows.Item(256).Delete()
